$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: bump the "Förändrad" (changed) date in column C for every data row
# (row 2 through the last used row) from 2023-09-14 (serial 45183) to
# 2023-09-15 (serial 45184) - one day later.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45184
}

# Step 2: rewrite the HYPERLINK formulas in columns S, T, V, W, X, Y for rows 2-5
# to add a friendly-text second argument (and promote column Y from a literal
# inline string to a real formula). These exact strings (including the
# mis-quoted column-S formula) mirror the source data verbatim.

# Row 2
$ws.Cells.Item(2, 19).Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_OSTERSUND/artfynd/A 30683-2023.xlsx, "A 30683-2023"")'
$ws.Cells.Item(2, 20).Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_OSTERSUND/kartor/A 30683-2023.png", "A 30683-2023")'
$ws.Cells.Item(2, 22).Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_OSTERSUND/klagomål/A 30683-2023.docx", "A 30683-2023")'
$ws.Cells.Item(2, 23).Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_OSTERSUND/klagomålsmail/A 30683-2023.docx", "A 30683-2023")'
$ws.Cells.Item(2, 24).Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_OSTERSUND/tillsyn/A 30683-2023.docx", "A 30683-2023")'
$ws.Cells.Item(2, 25).Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_OSTERSUND/tillsynsmail/A 30683-2023.docx", "A 30683-2023")'

# Row 3
$ws.Cells.Item(3, 19).Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_OSTERSUND/artfynd/A 32699-2023.xlsx, "A 32699-2023"")'
$ws.Cells.Item(3, 20).Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_OSTERSUND/kartor/A 32699-2023.png", "A 32699-2023")'
$ws.Cells.Item(3, 22).Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_OSTERSUND/klagomål/A 32699-2023.docx", "A 32699-2023")'
$ws.Cells.Item(3, 23).Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_OSTERSUND/klagomålsmail/A 32699-2023.docx", "A 32699-2023")'
$ws.Cells.Item(3, 24).Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_OSTERSUND/tillsyn/A 32699-2023.docx", "A 32699-2023")'
$ws.Cells.Item(3, 25).Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_OSTERSUND/tillsynsmail/A 32699-2023.docx", "A 32699-2023")'

# Row 4
$ws.Cells.Item(4, 19).Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_OSTERSUND/artfynd/A 29992-2023.xlsx, "A 29992-2023"")'
$ws.Cells.Item(4, 20).Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_OSTERSUND/kartor/A 29992-2023.png", "A 29992-2023")'
$ws.Cells.Item(4, 22).Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_OSTERSUND/klagomål/A 29992-2023.docx", "A 29992-2023")'
$ws.Cells.Item(4, 23).Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_OSTERSUND/klagomålsmail/A 29992-2023.docx", "A 29992-2023")'
$ws.Cells.Item(4, 24).Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_OSTERSUND/tillsyn/A 29992-2023.docx", "A 29992-2023")'
$ws.Cells.Item(4, 25).Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_OSTERSUND/tillsynsmail/A 29992-2023.docx", "A 29992-2023")'

# Row 5
$ws.Cells.Item(5, 19).Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_BRACKE/artfynd/A 30834-2023.xlsx, "A 30834-2023"")'
$ws.Cells.Item(5, 20).Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_BRACKE/kartor/A 30834-2023.png", "A 30834-2023")'
$ws.Cells.Item(5, 22).Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_BRACKE/klagomål/A 30834-2023.docx", "A 30834-2023")'
$ws.Cells.Item(5, 23).Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_BRACKE/klagomålsmail/A 30834-2023.docx", "A 30834-2023")'
$ws.Cells.Item(5, 24).Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_BRACKE/tillsyn/A 30834-2023.docx", "A 30834-2023")'
$ws.Cells.Item(5, 25).Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_BRACKE/tillsynsmail/A 30834-2023.docx", "A 30834-2023")'
